$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.733.59'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '1.595.33'
$ws.Range('E3').Value = '  -1.93%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''208.44'
$ws.Range('E5').Value = '  -1.49%  '
$ws.Range('E6').Value = '  -2.66%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '''22.37'
$ws.Range('E8').Value = '  -3.29%  '
$ws.Range('E9').Value = '  -1.87%  '
$ws.Range('D10').Value = '''0.0594'
$ws.Range('E10').Value = '  -2.00%  '
$ws.Range('E11').Value = '  -1.56%  '
$ws.Range('D12').Value = '1.820.96'
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('D13').Value = '1.603.36'
$ws.Range('E13').Value = '  -1.49%  '
$ws.Range('E14').Value = '  -3.49%  '
$ws.Range('D15').Value = '''0.534'
$ws.Range('E15').Value = '  -3.70%  '
$ws.Range('D16').Value = '27.734.37'
$ws.Range('E16').Value = '  -0.67%  '
$ws.Range('D17').Value = '''63.53'
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('D18').Value = '''220.16'
$ws.Range('E18').Value = '  -3.33%  '
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('D20').Value = '''7.38'
$ws.Range('E20').Value = '  -3.03%  '
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('E22').Value = '  -4.08%  '
$ws.Range('D23').Value = '''9.74'
$ws.Range('E23').Value = '  -2.15%  '
$ws.Range('D24').Value = '''1.98'
$ws.Range('D25').Value = '''153.95'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').Value = '''6.79'
$ws.Range('E26').Value = '  -1.74%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = '''15.17'
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('D29').Value = '''0.105'
$ws.Range('E29').Value = '  -4.79%  '
$ws.Range('D30').Value = '''1.17'
$ws.Range('E30').Value = '  -1.30%  '
$ws.Range('E31').Value = '  -1.70%  '
$ws.Range('E32').Value = '  -4.92%  '
$ws.Range('D33').Value = '1.376.59'
$ws.Range('E33').Value = '  -2.81%  '
$ws.Range('D34').Value = '''2.97'
$ws.Range('E34').Value = '  -4.25%  '
$ws.Range('E35').Value = '  -3.98%  '
$ws.Range('D36').Value = '''0.972'
$ws.Range('E36').Value = '  -2.73%  '
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('E39').Value = '  -2.94%  '
$ws.Range('D40').Value = '''0.829'
$ws.Range('E40').Value = '  -2.45%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('E42').Value = '  -3.30%  '
$ws.Range('D43').Value = '''64.62'
$ws.Range('E43').Value = '  -1.35%  '
$ws.Range('E44').Value = '  +2.57%  '
$ws.Range('E45').Value = '  -3.61%  '
$ws.Range('D46').Value = '''1.74'
$ws.Range('E46').Value = '  -4.88%  '
$ws.Range('D47').Value = '1.731.98'
$ws.Range('E47').Value = '  -2.11%  '
$ws.Range('D48').Value = '''87.03'
$ws.Range('E48').Value = '  -1.94%  '
$ws.Range('D49').Value = '0.0₆0102'
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('D50').Value = '''0.0966'
$ws.Range('E50').Value = '  -3.99%  '
$ws.Range('E51').Value = '  -1.25%  '
